$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update CP001 row (row 2): E2 "h5" -> "No es posible ingresar"
$ws.Range("E2").Value = "No es posible ingresar"

# --- Update CP002 row (row 3): C3 "div" -> new "no results" message
$ws.Range("C3").Value = 'No encontramos resultados para tu consulta. Te recomendamos usar frases cortas y palabras claves. Ej.: "caja de ahorro".'

# --- Row 4 (CP003_Eminent): clear leftover C4/E4 data that no longer applies
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()

# --- Row 5: replace the old placeholder CP004 test data with the new
#     CP004_Sucursal006 case
$ws.Range("A5").Value = "CP004_Sucursal006"
$ws.Range("B5").Value = "FLORES"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = "RIVADAVIA"
$ws.Range("E5").Value = "AV. RIVADAVIA 7121 "

# --- Rows 6-9: drop the old placeholder Dato1/2/3 CPx values, keep only
#     the CP00x label in column A
$ws.Range("B6:E6").ClearContents()
$ws.Range("B7:E7").ClearContents()
$ws.Range("B8:E8").ClearContents()
$ws.Range("B9:E9").ClearContents()

# --- Widen column D to fit the new "RIVADAVIA" data
$ws.Columns("D").ColumnWidth = 19.6

# --- Update the active selection to C6, matching the author's last position
$ws.Range("C6").Select()
